# Update odds values in Sheet1 to match the latest FlashScore export.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 (Tigre - Banfield)
$ws.Range("G6").Value = 2.3
$ws.Range("H6").Value = 2.9
$ws.Range("I6").Value = 3.5
$ws.Range("J6").Value = 3.2
$ws.Range("K6").Value = 1.83
$ws.Range("L6").Value = 4.5
$ws.Range("W6").Value = 2.38
$ws.Range("X6").Value = 1.53
$ws.Range("Y6").Value = 5.5
$ws.Range("Z6").Value = 9
$ws.Range("AB6").Value = 21
$ws.Range("AC6").Value = 26
$ws.Range("AG6").Value = 21
$ws.Range("AJ6").Value = 7
$ws.Range("AK6").Value = 15

# Row 8 (Racing Cordoba - Alvarado)
$ws.Range("G8").Value = 1.83
$ws.Range("H8").Value = 3.1
$ws.Range("J8").Value = 2.63
$ws.Range("K8").Value = 1.83
$ws.Range("AH8").Value = 126
$ws.Range("AJ8").Value = 8.5

# Row 11 (Colo Colo - O'Higgins)
$ws.Range("M11").Value = 1.03
$ws.Range("O11").Value = 1.19
$ws.Range("P11").Value = 4
$ws.Range("S11").Value = 2.75
$ws.Range("T11").Value = 1.37
$ws.Range("AB11").Value = 11
$ws.Range("AI11").Value = 251
$ws.Range("AN11").Value = 41

# Row 12 (Junior - Envigado)
$ws.Range("Q12").Value = 1.93
$ws.Range("R12").Value = 1.93
$ws.Range("S12").Value = 3.25
$ws.Range("T12").Value = 1.33
$ws.Range("AP12").Value = 1.46
$ws.Range("AQ12").Value = 2.7
$ws.Range("AR12").Value = 2.6
$ws.Range("AS12").Value = 1.48

# Row 17 (Tecnico U. - Libertad)
$ws.Range("G17").Value = 1.62
$ws.Range("H17").Value = 3.75
$ws.Range("W17").Value = 2.1
$ws.Range("X17").Value = 1.67
$ws.Range("Z17").Value = 7
$ws.Range("AF17").Value = 7.5
$ws.Range("AG17").Value = 21
$ws.Range("AK17").Value = 26
$ws.Range("AO17").Value = 51

# Row 26 (Sportivo Trinidense - General Caballero JLM)
$ws.Range("G26").Value = 2.2
$ws.Range("H26").Value = 3.3
$ws.Range("I26").Value = 3.3
$ws.Range("J26").Value = 3
$ws.Range("Z26").Value = 10
$ws.Range("AA26").Value = 9.5
$ws.Range("AB26").Value = 21
$ws.Range("AD26").Value = 34
$ws.Range("AE26").Value = 8.5
$ws.Range("AG26").Value = 15
$ws.Range("AI26").Value = 351
$ws.Range("AK26").Value = 15
$ws.Range("AM26").Value = 34

# Row 27 (Binacional - Cusco)
$ws.Range("G27").Value = 2.45
$ws.Range("I27").Value = 2.75
$ws.Range("J27").Value = 3.1
$ws.Range("L27").Value = 3.4
$ws.Range("M27").Value = 1.04
$ws.Range("N27").Value = 10
$ws.Range("O27").Value = 1.25
$ws.Range("T27").Value = 1.27
$ws.Range("Y27").Value = 8.5
$ws.Range("Z27").Value = 12
$ws.Range("AA27").Value = 10
$ws.Range("AB27").Value = 23
$ws.Range("AC27").Value = 21
$ws.Range("AH27").Value = 41
$ws.Range("AJ27").Value = 9.5
$ws.Range("AL27").Value = 11
$ws.Range("AM27").Value = 29
$ws.Range("AN27").Value = 21
$ws.Range("AO27").Value = 29
